# Fruta / hortaliza, semanal
# Insert 3 new weekly rows of "Ciruela" (plum) price data at row 197,
# pushing the existing rows 197-204 down to 200-207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 197 (each Insert() pushes
# row 197 and below down by one, so calling it three times in a row
# yields three new blank rows at 197, 198, 199).
$ws.Rows.Item(197).Insert()
$ws.Rows.Item(197).Insert()
$ws.Rows.Item(197).Insert()

# New row 197: Angeleno / Especial
$ws.Cells.Item(197, 1).Value = 9
$ws.Cells.Item(197, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(197, 3).Value = "Metropolitana"
$ws.Cells.Item(197, 4).Value = 45021
$ws.Cells.Item(197, 5).Value = 13
$ws.Cells.Item(197, 6).Value = "Fruta"
$ws.Cells.Item(197, 7).Value = 100103
$ws.Cells.Item(197, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(197, 9).Value = 100103002
$ws.Cells.Item(197, 10).Value = "Ciruela"
$ws.Cells.Item(197, 11).Value = "Angeleno"
$ws.Cells.Item(197, 12).Value = "Especial"
$ws.Cells.Item(197, 13).Value = 250
$ws.Cells.Item(197, 14).Value = 10500
$ws.Cells.Item(197, 15).Value = 10500
$ws.Cells.Item(197, 16).Value = 10500
$ws.Cells.Item(197, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(197, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(197, 19).Value = 700
$ws.Cells.Item(197, 20).Value = 15

# New row 198: Angeleno / Primera
$ws.Cells.Item(198, 1).Value = 9
$ws.Cells.Item(198, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 45021
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = "Fruta"
$ws.Cells.Item(198, 7).Value = 100103
$ws.Cells.Item(198, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(198, 9).Value = 100103002
$ws.Cells.Item(198, 10).Value = "Ciruela"
$ws.Cells.Item(198, 11).Value = "Angeleno"
$ws.Cells.Item(198, 12).Value = "Primera"
$ws.Cells.Item(198, 13).Value = 290
$ws.Cells.Item(198, 14).Value = 9000
$ws.Cells.Item(198, 15).Value = 9000
$ws.Cells.Item(198, 16).Value = 9000
$ws.Cells.Item(198, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(198, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(198, 19).Value = 600
$ws.Cells.Item(198, 20).Value = 15

# New row 199: Angeleno / Segunda
$ws.Cells.Item(199, 1).Value = 9
$ws.Cells.Item(199, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(199, 3).Value = "Metropolitana"
$ws.Cells.Item(199, 4).Value = 45021
$ws.Cells.Item(199, 5).Value = 13
$ws.Cells.Item(199, 6).Value = "Fruta"
$ws.Cells.Item(199, 7).Value = 100103
$ws.Cells.Item(199, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(199, 9).Value = 100103002
$ws.Cells.Item(199, 10).Value = "Ciruela"
$ws.Cells.Item(199, 11).Value = "Angeleno"
$ws.Cells.Item(199, 12).Value = "Segunda"
$ws.Cells.Item(199, 13).Value = 200
$ws.Cells.Item(199, 14).Value = 7500
$ws.Cells.Item(199, 15).Value = 7500
$ws.Cells.Item(199, 16).Value = 7500
$ws.Cells.Item(199, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(199, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(199, 19).Value = 500
$ws.Cells.Item(199, 20).Value = 15
